{"js": "// Update the worksheet date and the 25 three-digit \u00f7 one-digit problems\n// to the new day's values. Every \"old\" string below is unique in the\n// document, so a body-wide search/replace per pair is unambiguous.\nconst body = context.document.body;\nconst replacements = [\n  [\"2024-03-08 Friday\", \"2024-03-09 Saturday\"],\n  [\"394\u00f78=\", \"119\u00f79=\"],\n  [\"698\u00f79=\", \"571\u00f75=\"],\n  [\"207\u00f79=\", \"873\u00f74=\"],\n  [\"124\u00f75=\", \"779\u00f73=\"],\n  [\"124\u00f79=\", \"139\u00f76=\"],\n  [\"688\u00f79=\", \"508\u00f77=\"],\n  [\"192\u00f72=\", \"837\u00f77=\"],\n  [\"654\u00f76=\", \"156\u00f78=\"],\n  [\"788\u00f72=\", \"173\u00f72=\"],\n  [\"647\u00f76=\", \"683\u00f73=\"],\n  [\"137\u00f74=\", \"565\u00f75=\"],\n  [\"238\u00f78=\", \"365\u00f76=\"],\n  [\"663\u00f75=\", \"396\u00f78=\"],\n  [\"556\u00f75=\", \"582\u00f77=\"],\n  [\"165\u00f76=\", \"176\u00f76=\"],\n  [\"478\u00f77=\", \"998\u00f72=\"],\n  [\"700\u00f79=\", \"653\u00f72=\"],\n  [\"473\u00f74=\", \"578\u00f77=\"],\n  [\"793\u00f72=\", \"905\u00f72=\"],\n  [\"125\u00f75=\", \"746\u00f79=\"],\n  [\"223\u00f74=\", \"358\u00f76=\"],\n  [\"915\u00f77=\", \"200\u00f74=\"],\n  [\"742\u00f78=\", \"612\u00f78=\"],\n  [\"723\u00f76=\", \"358\u00f78=\"],\n  [\"196\u00f74=\", \"583\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 three-digit \u00f7 one-digit problems\n# to the new day's values. Every \"old\" string below is unique in the\n# document, so a document-wide Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2024-03-08 Friday', '2024-03-09 Saturday'),\n    @('394\u00f78=', '119\u00f79='),\n    @('698\u00f79=', '571\u00f75='),\n    @('207\u00f79=', '873\u00f74='),\n    @('124\u00f75=', '779\u00f73='),\n    @('124\u00f79=', '139\u00f76='),\n    @('688\u00f79=', '508\u00f77='),\n    @('192\u00f72=', '837\u00f77='),\n    @('654\u00f76=', '156\u00f78='),\n    @('788\u00f72=', '173\u00f72='),\n    @('647\u00f76=', '683\u00f73='),\n    @('137\u00f74=', '565\u00f75='),\n    @('238\u00f78=', '365\u00f76='),\n    @('663\u00f75=', '396\u00f78='),\n    @('556\u00f75=', '582\u00f77='),\n    @('165\u00f76=', '176\u00f76='),\n    @('478\u00f77=', '998\u00f72='),\n    @('700\u00f79=', '653\u00f72='),\n    @('473\u00f74=', '578\u00f77='),\n    @('793\u00f72=', '905\u00f72='),\n    @('125\u00f75=', '746\u00f79='),\n    @('223\u00f74=', '358\u00f76='),\n    @('915\u00f77=', '200\u00f74='),\n    @('742\u00f78=', '612\u00f78='),\n    @('723\u00f76=', '358\u00f78='),\n    @('196\u00f74=', '583\u00f77=')\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n}\n"}
